$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Captured_Values")

for ($r = 21; $r -le 29; $r++) {
    $ws.Cells.Item($r, 1).Value = 123456789
    $ws.Cells.Item($r, 2).Value = "Real Programmers Count 0123456789 From Zero"
}
